# Update "想去人数" (attendance) figures across sheets, per upstream data refresh.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheets index 1) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F6").Value  = 666
$wsExhibit.Range("F10").Value = 164
$wsExhibit.Range("F12").Value = 7854
$wsExhibit.Range("F18").Value = 270
$wsExhibit.Range("F21").Value = 279
$wsExhibit.Range("F22").Value = 9680
$wsExhibit.Range("F24").Value = 274
$wsExhibit.Range("F38").Value = 816
$wsExhibit.Range("F39").Value = 3987
$wsExhibit.Range("F41").Value = 1082
$wsExhibit.Range("F48").Value = 49

# --- Sheet "演出" (Worksheets index 2) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 3

# --- Sheet "全部类型" (Worksheets index 4) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F7").Value  = 666
$wsAll.Range("F13").Value = 164
$wsAll.Range("F15").Value = 7854
$wsAll.Range("F19").Value = 270
$wsAll.Range("F21").Value = 279
$wsAll.Range("F22").Value = 9680
$wsAll.Range("F24").Value = 274
$wsAll.Range("F35").Value = 816
$wsAll.Range("F37").Value = 3987
$wsAll.Range("F39").Value = 1082
$wsAll.Range("F48").Value = 49

$wb.Save()
